$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold textual price values (e.g. "1.00", "67.27").
# Excel auto-converts numeric-looking strings assigned via .Value into
# real numbers, which would lose formatting such as trailing zeros or
# thousands separators. Force these specific cells to Text format first
# so the assigned strings are preserved exactly, matching the source data.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.229.65"
$ws.Range("E2").Value = "  -0.96%  "
$ws.Range("D3").Value = "2.611.24"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "592.02"
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("D6").Value = "151.99"
$ws.Range("E6").Value = "  -2.24%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "0.553"
$ws.Range("E8").Value = "  +0.68%  "
$ws.Range("D9").Value = "2.609.86"
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("D10").Value = "0.122"
$ws.Range("E10").Value = "  -2.87%  "
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("D12").Value = "5.19"
$ws.Range("E12").Value = "  -1.37%  "
$ws.Range("D13").Value = "0.347"
$ws.Range("E13").Value = "  -2.53%  "
$ws.Range("D14").Value = "27.58"
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("D15").Value = "3.081.01"
$ws.Range("E15").Value = "  -0.35%  "
$ws.Range("E16").Value = "  -4.18%  "
$ws.Range("D17").Value = "67.050.50"
$ws.Range("E17").Value = "  -1.02%  "
$ws.Range("D18").Value = "2.607.76"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").Value = "366.47"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("D20").Value = "11.03"
$ws.Range("E20").Value = "  -1.41%  "
$ws.Range("D21").Value = "7.37"
$ws.Range("E21").Value = "  -4.23%  "
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("E23").Value = "  +2.40%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("B25").Value = "Aptos"
$ws.Range("C25").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D25").Value = "10.10"
$ws.Range("E25").Value = "  +2.24%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").Value = "67.27"
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("D27").Value = "2.740.34"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("B28").Value = "Bittensor"
$ws.Range("C28").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D28").Value = "584.66"
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "0.0000101"
$ws.Range("E30").Value = "  -3.35%  "
$ws.Range("D31").Value = "1.38"
$ws.Range("E31").Value = "  -3.11%  "
$ws.Range("D32").Value = "7.70"
$ws.Range("E32").Value = "  -3.43%  "
$ws.Range("E33").Value = "  -3.25%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").Value = "0.122"
$ws.Range("E35").Value = "  -8.50%  "
$ws.Range("D36").Value = "1.50"
$ws.Range("E36").Value = "  -1.32%  "
$ws.Range("D37").Value = "4.84"
$ws.Range("E37").Value = "  -2.12%  "
$ws.Range("D38").Value = "156.12"
$ws.Range("E38").Value = "  +0.49%  "
$ws.Range("D39").Value = "18.89"
$ws.Range("E39").Value = "  -2.59%  "
$ws.Range("E40").Value = "  -1.12%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "1.85"
$ws.Range("E41").Value = "  -1.25%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").Value = "5.25"
$ws.Range("E42").Value = "  -2.06%  "
$ws.Range("D43").Value = "2.58"
$ws.Range("E43").Value = "  +1.13%  "
$ws.Range("D44").Value = "16.79"
$ws.Range("E44").Value = "  +2.11%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "155.09"
$ws.Range("E46").Value = "  -1.15%  "
$ws.Range("D47").Value = "0.0₆0294"
$ws.Range("E47").Value = "  -0.88%  "
$ws.Range("D48").Value = "3.72"
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("D49").Value = "21.84"
$ws.Range("E49").Value = "  +4.89%  "
$ws.Range("D50").Value = "1.70"
$ws.Range("E50").Value = "  -1.74%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.0785"
$ws.Range("E51").Value = "  +0.49%  "
